$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.047.54"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.808.69"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "'232.43"
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").Value = "'40.23"
$ws.Range("D9").Value = "'0.324"
$ws.Range("E9").Value = "  +5.62%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "'0.0997"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "2.071.78"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "1.822.21"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "'0.665"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'11.03"
$ws.Range("E15").Value = "  -5.36%  "
$ws.Range("D16").Value = "'4.65"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("D17").Value = "35.017.53"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "'69.68"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "0.0₃0790"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'237.60"
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("D21").Value = "'11.91"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").Value = "'4.72"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("D25").Value = "'171.99"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "'7.85"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").Value = "'0.121"
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("E29").Value = "  +19.31%  "
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").Value = "'4.14"
$ws.Range("E31").Value = "  +5.36%  "
$ws.Range("D32").Value = "'0.0555"
$ws.Range("E32").Value = "  +4.35%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("E34").Value = "  -5.97%  "
$ws.Range("E35").Value = "  +5.42%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.685"
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("B37").Value = "Aave"
$ws.Range("C37").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D37").Value = "'91.78"
$ws.Range("E37").Value = "  +3.29%  "
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").Value = "1.311.54"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'2.29"
$ws.Range("E43").Value = "  -4.62%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'14.37"
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("D46").Value = "'6.32"
$ws.Range("E46").Value = "  +5.49%  "
$ws.Range("D47").Value = "'0.0513"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").Value = "1.987.91"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "'0.0658"
$ws.Range("E50").Value = "  +5.65%  "
$ws.Range("D51").Value = "'99.39"
$ws.Range("E51").Value = "  -4.60%  "

foreach ($addr in @("D5","D8","D9","D11","D14","D15","D16","D18","D20","D21","D22","D25","D26","D28","D31","D32","D36","D37","D43","D44","D46","D47","D50","D51")) { $ws.Range($addr).Style = "Normal" }
